# Practica 1: se agrego deteccion de anomalias
#
# - Renames the headers on "normales" to "Estatura(metros)" / "Edad(años)"
# - Carves out a validation slice (last 6 rows of "normales" plus two new
#   synthetic points) into a new sheet "valtest(normales)" placed right
#   after "normales"
# - Leaves "anomalias" sheet as-is (still "Estatura" / "Edad" headers)

$wb = $excel.ActiveWorkbook

$normales = $wb.Worksheets.Item("normales")

# --- 1. Update the "normales" header text -------------------------------
$normales.Range("A1").Value = "Estatura(metros)"
$normales.Range("B1").Value = "Edad(años)"

# --- 2. Insert the new "valtest(normales)" sheet right after "normales" --
$valtest = $wb.Worksheets.Add($null, $normales)
$valtest.Name = "valtest(normales)"

$valtest.Range("A1").Value = "Estatura(metros)"
$valtest.Range("B1").Value = "Edad(años)"

$valtestData = @(
    @("1.82", 25),
    @("1.80", 27),
    @("1.60", 31),
    @("1.60", 35),
    @("1.82", 30),
    @("1.76", 32),
    @("1.79", 31),
    @("1.61", 31)
)

$r = 2
foreach ($row in $valtestData) {
    $valtest.Cells.Item($r, 1).Value = "'" + $row[0]
    $valtest.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# --- 3. Remove the rows that were moved out of "normales" (31-36) -------
$normales.Range("A31:B36").Clear()
